# Refresh the crypto price/volume table (cols D: Price, E: Volume(1h))
# on Sheet1 with the latest scraped values. Rows 12/13 also swap their
# Coin/Link/Price/Volume because TRON overtook Polkadot in the ranking.
#
# Several "Price" values look numeric (e.g. 74.25) but must stay stored
# as plain text, matching every other cell in this column (prices with
# thousands separators like "42.201.87" can never be parsed as numbers,
# but simple decimals would be auto-converted unless we force a Text
# format first). We flip the cell to Text, assign the literal string,
# then restore the "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.201.87"
$ws.Range("E2").Value = "  -0.59%  "

$ws.Range("D3").Value = "2.235.54"
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.63%  "

$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.25"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  -2.23%  "

$ws.Range("E10").Value = "  -2.02%  "

$ws.Range("E11").Value = "  +1.35%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.103"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.40%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.44%  "

$ws.Range("D14").Value = "2.573.22"
$ws.Range("E14").Value = "  +0.04%  "

$ws.Range("E15").Value = "  -0.53%  "

$ws.Range("E16").Value = "  -1.49%  "

$ws.Range("D17").Value = "2.229.11"
$ws.Range("E17").Value = "  -0.75%  "

$ws.Range("D18").Value = "42.070.31"
$ws.Range("E18").Value = "  -0.30%  "

$ws.Range("E19").Value = "  +1.67%  "

$ws.Range("E20").Value = "  +1.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +12.94%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("E24").Value = "  -6.22%  "

$ws.Range("E25").Value = "  -0.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.78%  "

$ws.Range("E27").Value = "  -1.05%  "

$ws.Range("E28").Value = "  -0.78%  "

$ws.Range("E29").Value = "  -2.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.49%  "

$ws.Range("E31").Value = "  -1.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0805"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.25%  "

$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("E36").Value = "  -6.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.81%  "

$ws.Range("E38").Value = "  -1.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.28%  "

$ws.Range("E40").Value = "  -1.16%  "

$ws.Range("E41").Value = "  +0.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.45%  "

$ws.Range("E43").Value = "  +0.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.92%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.37%  "

$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("E48").Value = "  -1.51%  "

$ws.Range("E49").Value = "  +0.12%  "

$ws.Range("E50").Value = "  -1.60%  "

$ws.Range("D51").Value = "2.443.94"
$ws.Range("E51").Value = "  -0.06%  "
